$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2023-10-08 10:42:28"
$dataset = "hatespeech"
$setting = "def"

$rows = @(
    @{ Row = 38; Dim = 10; M1 = 0.6570267269092438 },
    @{ Row = 39; Dim = 40; M1 = 0.3301035082478216 },
    @{ Row = 40; Dim = 30; M1 = 0.4054775467285684 },
    @{ Row = 41; Dim = 20; M1 = 0.5007948904606515 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $timestamp
    $ws.Cells.Item($rowNum, 2).Value = $dataset
    $ws.Cells.Item($rowNum, 3).Value = $setting
    $ws.Cells.Item($rowNum, 4).Value = $r.Dim
    $ws.Cells.Item($rowNum, 5).Value = $r.M1
}
